$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 10-13 (shrinks dimension to A1:T9)
$ws.Range("A10:T13").Delete() | Out-Null

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Cadm1"
$ws.Range("C2").Value = "Crtam"
$ws.Range("D2").Value = "MuSCs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 2.152962333333333
$ws.Range("H2").Value = 6.458887000000001
$ws.Range("I2").Value = 0.1024890697041326
$ws.Range("J2").Value = 0.1024890697041326
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.05032833333333334
$ws.Range("N2").Value = 0.150985
$ws.Range("O2").Value = 0.7147658790836833
$ws.Range("P2").Value = 0.7147658790836833
$ws.Range("Q2").Value = 0.1083550059661111
$ws.Range("R2").Value = 0.9751950536950001
$ws.Range("S2").Value = 0.07325569000354321
$ws.Range("T2").Value = 0.07325569000354322

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Cadm1"
$ws.Range("C3").Value = "Crtam"
$ws.Range("D3").Value = "Resolving-Mac"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 2.152962333333333
$ws.Range("H3").Value = 6.458887000000001
$ws.Range("I3").Value = 0.1024890697041326
$ws.Range("J3").Value = 0.1024890697041326
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.020084
$ws.Range("N3").Value = 0.060252
$ws.Range("O3").Value = 0.2852341209163167
$ws.Range("P3").Value = 0.2852341209163167
$ws.Range("Q3").Value = 0.04324009550266667
$ws.Range("R3").Value = 0.389160859524
$ws.Range("S3").Value = 0.02923337970058936
$ws.Range("T3").Value = 0.02923337970058937

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Cadm1"
$ws.Range("C4").Value = "Crtam"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.004511666666666666
$ws.Range("H4").Value = 0.013535
$ws.Range("I4").Value = 0.0002147722290923241
$ws.Range("J4").Value = 0.0002147722290923242
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.05032833333333334
$ws.Range("N4").Value = 0.150985
$ws.Range("O4").Value = 0.7147658790836833
$ws.Range("P4").Value = 0.7147658790836833
$ws.Range("Q4").Value = 0.0002270646638888889
$ws.Range("R4").Value = 0.002043581975
$ws.Range("S4").Value = 0.0001535118611299373
$ws.Range("T4").Value = 0.0001535118611299373

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Cadm1"
$ws.Range("C5").Value = "Crtam"
$ws.Range("D5").Value = "Resolving-Mac"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.004511666666666666
$ws.Range("H5").Value = 0.013535
$ws.Range("I5").Value = 0.0002147722290923241
$ws.Range("J5").Value = 0.0002147722290923242
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.020084
$ws.Range("N5").Value = 0.060252
$ws.Range("O5").Value = 0.2852341209163167
$ws.Range("P5").Value = 0.2852341209163167
$ws.Range("Q5").Value = 0.00009061231333333334
$ws.Range("R5").Value = 0.0008155108200000001
$ws.Range("S5").Value = 0.00006126036796238686
$ws.Range("T5").Value = 0.00006126036796238687

# Row 6
$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "Cadm1"
$ws.Range("C6").Value = "Crtam"
$ws.Range("D6").Value = "MuSCs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 3.553611
$ws.Range("H6").Value = 10.660833
$ws.Range("I6").Value = 0.1691651915323982
$ws.Range("J6").Value = 0.1691651915323982
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.05032833333333334
$ws.Range("N6").Value = 0.150985
$ws.Range("O6").Value = 0.7147658790836833
$ws.Range("P6").Value = 0.7147658790836833
$ws.Range("Q6").Value = 0.178847318945
$ws.Range("R6").Value = 1.609625870505
$ws.Range("S6").Value = 0.1209135068360142
$ws.Range("T6").Value = 0.1209135068360142

# Row 7
$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Cadm1"
$ws.Range("C7").Value = "Crtam"
$ws.Range("D7").Value = "Resolving-Mac"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 3.553611
$ws.Range("H7").Value = 10.660833
$ws.Range("I7").Value = 0.1691651915323982
$ws.Range("J7").Value = 0.1691651915323982
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.020084
$ws.Range("N7").Value = 0.060252
$ws.Range("O7").Value = 0.2852341209163167
$ws.Range("P7").Value = 0.2852341209163167
$ws.Range("Q7").Value = 0.071370723324
$ws.Range("R7").Value = 0.642336509916
$ws.Range("S7").Value = 0.04825168469638394
$ws.Range("T7").Value = 0.04825168469638395

# Row 8
$ws.Range("A8").Value = "Resolving-Mac"
$ws.Range("B8").Value = "Cadm1"
$ws.Range("C8").Value = "Crtam"
$ws.Range("D8").Value = "MuSCs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 15.29566566666667
$ws.Range("H8").Value = 45.886997
$ws.Range("I8").Value = 0.7281309665343768
$ws.Range("J8").Value = 0.7281309665343769
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.05032833333333334
$ws.Range("N8").Value = 0.150985
$ws.Range("O8").Value = 0.7147658790836833
$ws.Range("P8").Value = 0.7147658790836833
$ws.Range("Q8").Value = 0.7698053602272222
$ws.Range("R8").Value = 6.928248242045001
$ws.Range("S8").Value = 0.5204431703829958
$ws.Range("T8").Value = 0.5204431703829959

# Row 9
$ws.Range("A9").Value = "Resolving-Mac"
$ws.Range("B9").Value = "Cadm1"
$ws.Range("C9").Value = "Crtam"
$ws.Range("D9").Value = "Resolving-Mac"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 15.29566566666667
$ws.Range("H9").Value = 45.886997
$ws.Range("I9").Value = 0.7281309665343768
$ws.Range("J9").Value = 0.7281309665343769
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.020084
$ws.Range("N9").Value = 0.060252
$ws.Range("O9").Value = 0.2852341209163167
$ws.Range("P9").Value = 0.2852341209163167
$ws.Range("Q9").Value = 0.3071981492493334
$ws.Range("R9").Value = 2.764783343244
$ws.Range("S9").Value = 0.207687796151381
$ws.Range("T9").Value = 0.207687796151381
